# ---------------------------------------------------------------------------
# This workbook originally has a single sheet "27..43" holding one table of
# variant-rescue data. The update:
#   1. Duplicates that sheet (data + formatting) into a new sheet named
#      "All Results" placed right after the original, adds a "Helix" column
#      to it, and corrects a handful of values / regenerated sequences.
#   2. Replaces the content of the original "27..43" sheet with a fresh,
#      differently-shaped summary table (the "opool" table) that also gets
#      the new "Helix" column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Duplicate the existing sheet to become "All Results" -----------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "All Results"

# Add the new "Helix" header in column I, matching the existing header style
$ws2.Range("H1").Copy()
$ws2.Range("I1").PasteSpecial(-4122)
$ws2.Range("I1").Value = "Helix"

# The stochastic helix sequences were regenerated; update the two affected
# "Sequence" values and the resulting "Structure" text
$ws2.Range("B6").Value = "CTGGCGGAGTACCACGG"
$ws2.Range("B7").Value = "TGCCGGGAGTACCGGCA"
$ws2.Range("C7").Value = "((((((.....))))))"

# Corrected score/edit-distance values
$ws2.Range("H6").Value = 11
$ws2.Range("E7").Value = 5.19
$ws2.Range("F7").Value = 0
$ws2.Range("G7").Value = 1
$ws2.Range("H7").Value = 8

# New "Helix" metadata column for every non-WT row
$helixVal = "[[0, 16], [1, 15], [2, 14], [3, 13], [4, 12], [5, 11]]"
$ws2.Range("I3").Value = $helixVal
$ws2.Range("I4").Value = $helixVal
$ws2.Range("I5").Value = $helixVal
$ws2.Range("I6").Value = $helixVal
$ws2.Range("I7").Value = $helixVal

# --- 2. Rebuild the original "27..43" sheet with the new opool table -------
# Grab the header style from the (still untouched) "All Results" sheet before
# wiping "27..43", then paste it across the new header row.
$ws2.Range("B1").Copy()
$ws1.Cells.Clear()
$ws1.Range("A1:H1").PasteSpecial(-4122)

$headers = @("Sequence", "Structure", "Variant Type", "Folding Energy", "Disruption Score", "Recovery Score", "Edit Distance", "Helix")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws1.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$rows = @(
    @("GGCGTCGAGTAGACGCC", ".................", "WT", 5.93, $null, $null, $null, $null),
    @("CTGCGGGAGTAGACGCC", "((((((.....))))))", "Left Flip", -0.17, 1, 0, 6, $helixVal),
    @("GGCGTCGAGTACCGCAG", "((((((.....))))))", "Right Flip", -0.17, 1, 0, 6, $helixVal),
    @("CTGCGGGAGTACCGCAG", ".................", "Flipped Rescue", 5.2, 0, 1, 12, $helixVal),
    @("CTGGCGGAGTACCACGG", ".................", "Stochastic Helix Disruption", -0.17, 1, 1, 11, $helixVal),
    @("TGCCGGGAGTACCGGCA", "((((((.....))))))", "Stochastic Helix Rescue", 5.19, 0, 1, 8, $helixVal)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($val -ne $null) {
            $ws1.Cells.Item($excelRow, $c + 1).Value = $val
        }
    }
}

# Keep the original sheet the active/selected one
$ws1.Activate()
